# Add a new bullet ("Commits on a single branch") right after the
# "Special point 'head' to the current branch" list item, inheriting
# that paragraph's list/indent formatting (ListParagraph style, numId 3).

$d = $word.ActiveDocument

# Locate the anchor paragraph by its text rather than assuming it is
# strictly the document's last paragraph.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Special point*current branch*") {
        $anchor = $p
    }
}

# Insert a new paragraph mark after the anchor; Word carries the
# anchor's paragraph formatting (style + numbering) onto the new one.
$anchor.Range.InsertParagraphAfter()

# Re-fetch the paragraph collection's last item (now the freshly
# inserted paragraph) and give it its text.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Commits on a single branch"
